# Update cryptos list (Price / Volume(1h) columns) with refreshed market data.
# For D-column prices that look like plain decimal numbers, force the cell to
# stay textual (NumberFormat "@") before assigning, so Excel doesn't silently
# convert the string into a floating point number (which would corrupt values
# like "19.10" -> 19.1 or introduce binary float rounding noise).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.538.16'
$ws.Range('E2').Value = '  +2.51%  '
$ws.Range('D3').Value = '2.704.22'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('E4').Value = '  +0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '527.04'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.94'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('E8').Value = '  +1.87%  '
$ws.Range('D9').Value = '2.730.19'
$ws.Range('E9').Value = '  +3.10%  '
$ws.Range('E10').Value = '  +6.31%  '
$ws.Range('E11').Value = '  +0.62%  '
$ws.Range('E12').Value = '  +0.54%  '
$ws.Range('E13').Value = '  +3.06%  '
$ws.Range('D14').Value = '3.182.61'
$ws.Range('E14').Value = '  +3.04%  '
$ws.Range('D15').Value = '60.548.78'
$ws.Range('E15').Value = '  +2.61%  '
$ws.Range('D16').Value = '2.858.79'
$ws.Range('E16').Value = '  +8.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '21.32'
$ws.Range('E17').Value = '  +1.38%  '
$ws.Range('E18').Value = '  +0.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '346.64'
$ws.Range('E19').Value = '  -0.80%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('E21').Value = '  +3.05%  '
$ws.Range('E22').Value = '  +4.53%  '
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.69'
$ws.Range('E24').Value = '  +3.28%  '
$ws.Range('E25').Value = '  +0.38%  '
$ws.Range('E26').Value = '  +4.31%  '
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').Value = '0.0₃0819'
$ws.Range('E28').Value = '  +1.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.29'
$ws.Range('E29').Value = '  +2.21%  '
$ws.Range('E30').Value = '  +8.64%  '
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.10'
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.12'
$ws.Range('E34').Value = '  +0.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.26'
$ws.Range('E35').Value = '  +6.01%  '
$ws.Range('E36').Value = '  +8.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.942'
$ws.Range('E37').Value = '  -3.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.873'
$ws.Range('E38').Value = '  +3.31%  '
$ws.Range('E39').Value = '  +7.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '37.11'
$ws.Range('E40').Value = '  +1.01%  '
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '283.34'
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.10'
$ws.Range('E43').Value = '  +2.36%  '
$ws.Range('E44').Value = '  +0.48%  '
$ws.Range('E45').Value = '  -0.21%  '
$ws.Range('E46').Value = '  +0.44%  '
$ws.Range('D47').Value = '2.141.46'
$ws.Range('E47').Value = '  +8.01%  '
$ws.Range('E48').Value = '  +2.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '4.81'
$ws.Range('E49').Value = '  +3.38%  '
$ws.Range('E50').Value = '  +1.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0233'
$ws.Range('E51').Value = '  +1.49%  '
